$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Minimum Nr 3005 eingeführt, damit Sortierung etwas anzeigt
$ws.Range("A11").Value = 3005

# Update selection to B14
$ws.Range("B14").Select()

# Adjust window size/position (best-effort; mirrors the author's resized window)
$win = $excel.ActiveWindow
$win.Width = 7950
$win.Height = 11385
$win.Left = -28500
$win.Top = 2580
